$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "hISZs664"
$ws.Range("B2").Value = 23082422
$ws.Range("C2").Value = "gswjxye50"
$ws.Range("D2").Value = "v%6Xx3F#"
$ws.Range("F2").Value = "dVhilrjt"
$ws.Range("G2").Value = "TeFY"
